# Update gh-pages to output generated at 456a3b4
# Refresh "想去人数" (interest count) and "最低票价" (min price) figures
# across the 展览 (Exhibition), 演出 (Performance) and 全部类型 (All Types)
# sheets. 全部类型 aggregates the other sheets' rows, so every event that
# changed gets updated in both its home sheet and the All Types sheet.

$wb = $excel.ActiveWorkbook

function Set-F {
    param($ws, [int]$row, [double]$value)
    $ws.Cells.Item($row, 6).Value = $value
}

function Set-G {
    param($ws, [int]$row, [double]$value)
    $ws.Cells.Item($row, 7).Value = $value
}

# --- 展览 (Exhibition) sheet ---
$wsExpo = $wb.Worksheets.Item("展览")
Set-F $wsExpo 2  1837
Set-F $wsExpo 3  406
Set-G $wsExpo 3  45
Set-F $wsExpo 4  1497
Set-F $wsExpo 6  383
Set-G $wsExpo 6  60
Set-F $wsExpo 8  13216
Set-F $wsExpo 9  13080
Set-F $wsExpo 11 768
Set-F $wsExpo 13 549
Set-F $wsExpo 15 648
Set-F $wsExpo 16 2068
Set-F $wsExpo 18 31
Set-F $wsExpo 19 46

# --- 演出 (Performance) sheet ---
$wsShow = $wb.Worksheets.Item("演出")
Set-F $wsShow 7  107
Set-F $wsShow 9  13

# --- 全部类型 (All Types) sheet ---
$wsAll = $wb.Worksheets.Item("全部类型")
Set-F $wsAll 3  1837
Set-F $wsAll 4  406
Set-G $wsAll 4  45
Set-F $wsAll 5  1497
Set-F $wsAll 7  383
Set-G $wsAll 7  60
Set-F $wsAll 10 13216
Set-F $wsAll 11 13080
Set-F $wsAll 13 768
Set-F $wsAll 15 549
Set-F $wsAll 17 648
Set-F $wsAll 20 2068
Set-F $wsAll 22 31
Set-F $wsAll 23 46
Set-F $wsAll 31 107
Set-F $wsAll 33 13
